$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 76 (existing rows 76-95 shift down to 78-97)
$ws.Rows.Item(76).Insert()
$ws.Rows.Item(76).Insert()

# Match the date-column number format used by the rest of the table
$ws.Range("D76").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D77").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 76 - new Chirimoya "Primera" entry
$ws.Range("A76").Value = 10
$ws.Range("B76").Value = "Vega Modelo de Temuco"
$ws.Range("C76").Value = "La Araucanía"
$ws.Range("D76").Value = 44508
$ws.Range("E76").Value = 9
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100107
$ws.Range("H76").Value = "Otros"
$ws.Range("I76").Value = 100107002
$ws.Range("J76").Value = "Chirimoya"
$ws.Range("K76").Value = "Cultivar IV Región"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 125
$ws.Range("N76").Value = 3000
$ws.Range("O76").Value = 3000
$ws.Range("P76").Value = 3000
$ws.Range("Q76").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R76").Value = "Provincia del Elquí"
$ws.Range("S76").Value = 3000
$ws.Range("T76").Value = 1

# Row 77 - new Chirimoya "Segunda" entry
$ws.Range("A77").Value = 10
$ws.Range("B77").Value = "Vega Modelo de Temuco"
$ws.Range("C77").Value = "La Araucanía"
$ws.Range("D77").Value = 44508
$ws.Range("E77").Value = 9
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100107
$ws.Range("H77").Value = "Otros"
$ws.Range("I77").Value = 100107002
$ws.Range("J77").Value = "Chirimoya"
$ws.Range("K77").Value = "Cultivar IV Región"
$ws.Range("L77").Value = "Segunda"
$ws.Range("M77").Value = 95
$ws.Range("N77").Value = 2500
$ws.Range("O77").Value = 2500
$ws.Range("P77").Value = 2500
$ws.Range("Q77").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R77").Value = "Provincia del Elquí"
$ws.Range("S77").Value = 2500
$ws.Range("T77").Value = 1
